# =====================================================================
# Feat: Add all advanced tracking features - Email, Eftikad, Liturgy,
#       Excel Export, Announcements
#
# This script applies the workbook edit described by the diff:
#   1. "Students" sheet: append a new student row (id 150, Fady Massoud).
#   2. "Attendance Records" sheet: insert a new block of attendance rows
#      for 2026-01-31 ahead of the existing 2026-01-25 rows, and append
#      the remaining roster's 2026-01-25 attendance rows at the end.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------
# 1) Students sheet -- append row 151 (student id 150, "Fady Massoud")
# ----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Students")

$ws1.Cells.Item(151, 1).Value = 150
$ws1.Cells.Item(151, 2).Value = "Fady "
$ws1.Cells.Item(151, 3).Value = 8
$ws1.Cells.Item(151, 4).Value = "M"
$ws1.Cells.Item(151, 5).Value = "Massoud"
$ws1.Cells.Item(151, 6).Value = "213456yu"

# Columns G (parent_phone) and H (dob) hold digit-only / date-look-alike
# text in this workbook -- force Text format first so Excel does not
# coerce them into a Number / date serial value.
$ws1.Cells.Item(151, 7).NumberFormat = "@"
$ws1.Cells.Item(151, 7).Value = "231423546"
$ws1.Cells.Item(151, 8).NumberFormat = "@"
$ws1.Cells.Item(151, 8).Value = "2026-01-30"

# address / comments / pictures / last_call are blank for this row
$ws1.Cells.Item(151, 9).NumberFormat = "@"
$ws1.Cells.Item(151, 9).Value = ""
$ws1.Cells.Item(151, 10).NumberFormat = "@"
$ws1.Cells.Item(151, 10).Value = ""
$ws1.Cells.Item(151, 11).NumberFormat = "@"
$ws1.Cells.Item(151, 11).Value = ""
$ws1.Cells.Item(151, 12).NumberFormat = "@"
$ws1.Cells.Item(151, 12).Value = ""

# ----------------------------------------------------------------
# 2) Attendance Records sheet
# ----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Attendance Records")

# Insert 14 blank rows above row 2. This pushes the 3 existing rows
# (the 2026-01-25 attendance for Bishoy Arsalyos, Cristiano Gerges and
# Danny Soliman) down to rows 16-18, making room for the new
# 2026-01-31 attendance block.
$ws2.Range("A2:A15").EntireRow.Insert()

$ws2.Cells.Item(2, 1).Value = "Christin Wasef"
$ws2.Cells.Item(2, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(2, 2).Value = "2026-01-31"
$ws2.Cells.Item(2, 3).Value = "present"

$ws2.Cells.Item(3, 1).Value = "Jonathan Seif"
$ws2.Cells.Item(3, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(3, 2).Value = "2026-01-31"
$ws2.Cells.Item(3, 3).Value = "present"

$ws2.Cells.Item(4, 1).Value = "Mira Malty"
$ws2.Cells.Item(4, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(4, 2).Value = "2026-01-31"
$ws2.Cells.Item(4, 3).Value = "present"

$ws2.Cells.Item(5, 1).Value = "Novear Mikhael "
$ws2.Cells.Item(5, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(5, 2).Value = "2026-01-31"
$ws2.Cells.Item(5, 3).Value = "absent"

$ws2.Cells.Item(6, 1).Value = "Philopateer Kaldas"
$ws2.Cells.Item(6, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(6, 2).Value = "2026-01-31"
$ws2.Cells.Item(6, 3).Value = "present"

$ws2.Cells.Item(7, 1).Value = "Salah Salib"
$ws2.Cells.Item(7, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(7, 2).Value = "2026-01-31"
$ws2.Cells.Item(7, 3).Value = "absent"

$ws2.Cells.Item(8, 1).Value = "Shady Aziz"
$ws2.Cells.Item(8, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(8, 2).Value = "2026-01-31"
$ws2.Cells.Item(8, 3).Value = "absent"

$ws2.Cells.Item(9, 1).Value = "Shenouda Hanna"
$ws2.Cells.Item(9, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(9, 2).Value = "2026-01-31"
$ws2.Cells.Item(9, 3).Value = "present"

$ws2.Cells.Item(10, 1).Value = "Shenouda Saeed"
$ws2.Cells.Item(10, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(10, 2).Value = "2026-01-31"
$ws2.Cells.Item(10, 3).Value = "present"

$ws2.Cells.Item(11, 1).Value = "Thomas Keliny"
$ws2.Cells.Item(11, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(11, 2).Value = "2026-01-31"
$ws2.Cells.Item(11, 3).Value = "absent"

$ws2.Cells.Item(12, 1).Value = "Thomas Masoued"
$ws2.Cells.Item(12, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(12, 2).Value = "2026-01-31"
$ws2.Cells.Item(12, 3).Value = "present"

$ws2.Cells.Item(13, 1).Value = "Torres Ibrahim "
$ws2.Cells.Item(13, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(13, 2).Value = "2026-01-31"
$ws2.Cells.Item(13, 3).Value = "present"

$ws2.Cells.Item(14, 1).Value = "Wanas Youns/Abdelshahid"
$ws2.Cells.Item(14, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(14, 2).Value = "2026-01-31"
$ws2.Cells.Item(14, 3).Value = "absent"

$ws2.Cells.Item(15, 1).Value = "Yustos Bostros"
$ws2.Cells.Item(15, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(15, 2).Value = "2026-01-31"
$ws2.Cells.Item(15, 3).Value = "absent"

$ws2.Cells.Item(16, 1).Value = "Bishoy Arsalyos"
$ws2.Cells.Item(16, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(16, 2).Value = "2026-01-25"
$ws2.Cells.Item(16, 3).Value = "present"

$ws2.Cells.Item(17, 1).Value = "Cristiano Gerges (Attia?)"
$ws2.Cells.Item(17, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(17, 2).Value = "2026-01-25"
$ws2.Cells.Item(17, 3).Value = "present"

$ws2.Cells.Item(18, 1).Value = "Danny Soliman"
$ws2.Cells.Item(18, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(18, 2).Value = "2026-01-25"
$ws2.Cells.Item(18, 3).Value = "present"

$ws2.Cells.Item(19, 1).Value = "Fady Reda"
$ws2.Cells.Item(19, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(19, 2).Value = "2026-01-25"
$ws2.Cells.Item(19, 3).Value = "present"

$ws2.Cells.Item(20, 1).Value = "Georgino Bebawy"
$ws2.Cells.Item(20, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(20, 2).Value = "2026-01-25"
$ws2.Cells.Item(20, 3).Value = "present"

$ws2.Cells.Item(21, 1).Value = "Ishak Kamel"
$ws2.Cells.Item(21, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(21, 2).Value = "2026-01-25"
$ws2.Cells.Item(21, 3).Value = "present"

$ws2.Cells.Item(22, 1).Value = "John Yehia"
$ws2.Cells.Item(22, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(22, 2).Value = "2026-01-25"
$ws2.Cells.Item(22, 3).Value = "present"

$ws2.Cells.Item(23, 1).Value = "Jotham Shenouda "
$ws2.Cells.Item(23, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(23, 2).Value = "2026-01-25"
$ws2.Cells.Item(23, 3).Value = "present"

$ws2.Cells.Item(24, 1).Value = "Joyce Zaki"
$ws2.Cells.Item(24, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(24, 2).Value = "2026-01-25"
$ws2.Cells.Item(24, 3).Value = "absent"

$ws2.Cells.Item(25, 1).Value = "Justin Fakoury"
$ws2.Cells.Item(25, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(25, 2).Value = "2026-01-25"
$ws2.Cells.Item(25, 3).Value = "present"

$ws2.Cells.Item(26, 1).Value = "Karas Eshak Abdelmalak"
$ws2.Cells.Item(26, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(26, 2).Value = "2026-01-25"
$ws2.Cells.Item(26, 3).Value = "present"

$ws2.Cells.Item(27, 1).Value = "Karas Fares"
$ws2.Cells.Item(27, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(27, 2).Value = "2026-01-25"
$ws2.Cells.Item(27, 3).Value = "present"

$ws2.Cells.Item(28, 1).Value = "Karas Monir"
$ws2.Cells.Item(28, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(28, 2).Value = "2026-01-25"
$ws2.Cells.Item(28, 3).Value = "present"

$ws2.Cells.Item(29, 1).Value = "Karas Moura"
$ws2.Cells.Item(29, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(29, 2).Value = "2026-01-25"
$ws2.Cells.Item(29, 3).Value = "present"

$ws2.Cells.Item(30, 1).Value = "Karas Shehata "
$ws2.Cells.Item(30, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(30, 2).Value = "2026-01-25"
$ws2.Cells.Item(30, 3).Value = "present"

$ws2.Cells.Item(31, 1).Value = "Kyrollos Soliman"
$ws2.Cells.Item(31, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(31, 2).Value = "2026-01-25"
$ws2.Cells.Item(31, 3).Value = "present"

$ws2.Cells.Item(32, 1).Value = "Madonna Girgis"
$ws2.Cells.Item(32, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(32, 2).Value = "2026-01-25"
$ws2.Cells.Item(32, 3).Value = "absent"

$ws2.Cells.Item(33, 1).Value = "Maria Attiya"
$ws2.Cells.Item(33, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(33, 2).Value = "2026-01-25"
$ws2.Cells.Item(33, 3).Value = "absent"

$ws2.Cells.Item(34, 1).Value = "Mariam Ibrahim"
$ws2.Cells.Item(34, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(34, 2).Value = "2026-01-25"
$ws2.Cells.Item(34, 3).Value = "absent"

$ws2.Cells.Item(35, 1).Value = "Marly Abdelshehed"
$ws2.Cells.Item(35, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(35, 2).Value = "2026-01-25"
$ws2.Cells.Item(35, 3).Value = "absent"

$ws2.Cells.Item(36, 1).Value = "Marseleno Mina"
$ws2.Cells.Item(36, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(36, 2).Value = "2026-01-25"
$ws2.Cells.Item(36, 3).Value = "present"

$ws2.Cells.Item(37, 1).Value = "Monica Ramsis"
$ws2.Cells.Item(37, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(37, 2).Value = "2026-01-25"
$ws2.Cells.Item(37, 3).Value = "absent"

$ws2.Cells.Item(38, 1).Value = "Sabrina Kamel "
$ws2.Cells.Item(38, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(38, 2).Value = "2026-01-25"
$ws2.Cells.Item(38, 3).Value = "absent"

$ws2.Cells.Item(39, 1).Value = "Sandy Ibrahim"
$ws2.Cells.Item(39, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(39, 2).Value = "2026-01-25"
$ws2.Cells.Item(39, 3).Value = "absent"

$ws2.Cells.Item(40, 1).Value = "Verena Abdelmalak"
$ws2.Cells.Item(40, 2).NumberFormat = "@"  # keep date as text, not a serial number
$ws2.Cells.Item(40, 2).Value = "2026-01-25"
$ws2.Cells.Item(40, 3).Value = "absent"
